$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at 128 (pushes existing rows 128..227 down to 129..228,
# carrying their formatting/styles with them).
$ws.Rows.Item(128).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(128, 1).Value2  = 11
$ws.Cells.Item(128, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(128, 3).Value2  = "Bíobío"
$ws.Cells.Item(128, 4).Value2  = 45072
$ws.Cells.Item(128, 5).Value2  = 8
$ws.Cells.Item(128, 6).Value2  = 100112032
$ws.Cells.Item(128, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(128, 8).Value2  = "Sin especificar"
$ws.Cells.Item(128, 9).Value2  = "Primera"
$ws.Cells.Item(128, 10).Value2 = 220
$ws.Cells.Item(128, 11).Value2 = 7500
$ws.Cells.Item(128, 12).Value2 = 8000
$ws.Cells.Item(128, 13).Value2 = 7773
$ws.Cells.Item(128, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(128, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(128, 16).Value2 = 130
$ws.Cells.Item(128, 17).Value2 = 60
$ws.Cells.Item(128, 18).Value2 = "Hortaliza"
